# Taxi Kick - Doesn't allow kick while carrying passenger
#
# This updates the "TODO / notes" column (D) across a few sheets of the
# Development.xlsx tracker workbook, reflecting work that got done
# (kick-while-carrying-passenger handling) and reorganizes a couple of
# "NamedPipe" -> "ConPass" notes. It also switches the active sheet from
# CenDLL to CenTaxi.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# CTDLL sheet
# ---------------------------------------------------------------
$ctdll = $wb.Worksheets.Item("CTDLL")

# New remarks added next to "Send new position (when enter new cell)"
# and "Shutdown application" rows.
$ctdll.Range("D13").Value = "Needs to be looked at"
$ctdll.Range("D14").Value = "Needs completion"

# Remark on "Shutdown application" row is no longer needed.
$ctdll.Range("D16").ClearContents()

$ctdll.Range("D8").Select() | Out-Null

# ---------------------------------------------------------------
# CenTaxi sheet
# ---------------------------------------------------------------
$centaxi = $wb.Worksheets.Item("CenTaxi")

# Remark on "Receive taxi service status updates" row resolved.
$centaxi.Range("D16").ClearContents()

# The "Notify taxi about passenger assigned" / "Notify passenger about
# taxi assigned" rows used to share a "NamedPipe" remark; now only the
# second row keeps a remark, updated to "ConPass".
$centaxi.Range("D26").ClearContents()
$centaxi.Range("D27").Value = "ConPass"

# CenTaxi becomes the active sheet/tab (was CenDLL).
$centaxi.Activate() | Out-Null
$excel.ActiveWindow.ScrollRow = 23
$excel.ActiveWindow.ScrollColumn = 1
$centaxi.Range("C46").Select() | Out-Null
